# Update cryptocurrency price/volume data (Price column D, Volume(1h) column E)
# on Sheet1. Values that look like plain decimal numbers (e.g. "1.89") are
# forced to text via NumberFormat "@" so they are stored verbatim (matching
# the source data, which uses a dot as a thousands separator in some rows
# and must remain text, not be normalized as a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.950.55"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.702.05"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.41"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.89"
$ws.Range("E6").Value = "  +8.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "656.25"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  +3.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "3.696.96"
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.40"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.207"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").Value = "  +6.32%  "
$ws.Range("D15").Value = "4.389.14"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").Value = "96.796.49"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.08"
$ws.Range("E18").Value = "  +4.62%  "
$ws.Range("D19").Value = "3.692.82"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("E20").Value = "  +6.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.95"
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.527"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "525.39"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.12"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.91"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.42"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "3.897.60"
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.66"
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("E34").Value = "  +14.69%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.83"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "661.53"
$ws.Range("E38").Value = "  +7.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.607"
$ws.Range("E39").Value = "  +6.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.06"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.22"
$ws.Range("E41").Value = "  +19.75%  "
$ws.Range("E42").Value = "  +5.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.971"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.25"
$ws.Range("E45").Value = "  +19.10%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.452"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0461"
$ws.Range("E48").Value = "  +4.53%  "
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.84"
$ws.Range("E50").Value = "  +3.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.62"
$ws.Range("E51").Value = "  -0.03%  "
